$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Relates To ID" value for the "Alice in Wonderland - Title sequence" row
$ws.Range("J3").Value = "BE_001, BE_002"

# Move active selection to J7, matching the saved view state
$ws.Range("J7").Select()
